# Apply "rekcja" (case-government) additions to Arkusz1 (sheet1)
# New German phrase / Polish translation pairs are appended starting at row 216
# (row 215 intentionally left blank, matching the author's edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("sich reißen um A", "zabijać się o [coś] (pot.)"),
    @("aus-gehen von D", "pochodzić od [czegoś]/mieć swoje źródło w [czymś]"),
    @("jdm erpressen mit D", "szantażować kogoś [czymś]"),
    @("heran-kommen an A", "zbliżać się do [czegoś]"),
    @("sich entsinnen an A", "przypominać sobie [kogoś/coś]"),
    @("jdn unterstützen bei D/in D", "pomagać komuś przy [czymś]/w [czymś]"),
    @("neigen zu D", "mieć tendencję do [czegoś]"),
    @("zurück-gehen auf A", "mieć początek w [czymś], sięgać [czegoś]"),
    @("staunen über A", "dziwić się [komuś/czemuś]"),
    @("sich widerspiegeln in D", "odbijać/odzwierciedlać się w [czymś]"),
    @("jdn aus-schließen aus D", "usunąć, wykluczyć z [czegoś]"),
    @("Zweifel haben an D", "mieć wątpliwości co do [czegoś]"),
    @("protestieren gegen A", "protestować przeciwko [czemuś]"),
    @("überreden zu D", "namówić do [czegoś]"),
    @("sich fernhalten von D", "trzymać się z daleka od [kogoś/czegoś]"),
    @("an-knüpfen an A", "przywiązać do [czegoś]"),
    @("basteln an D", "majsterkować przy [czymś]"),
    @("ein-steigen in A", "wsiąść do [pojazdu]"),
    @("hungrig sein auf A", "mieć apetyt na [coś]"),
    @("jdn ansprechen auf A", "zwrócić się do kogoś w sprawie [czegoś]"),
    @("kranken an D", "chorować na [coś]"),
    @("sich äußern zu D", "wyrazić swoje zdanie co do [czegoś]"),
    @("sich engagieren für A", "zaangażować się [w]"),
    @("sich erschrecken vor D", "przestraszyć się [kogoś/czegoś]"),
    @("sich zurechtfinden in D (z)", "orientować się w [mieście/otoczeniu]"),
    @("vorbei-schauen bei D", "zaglądać, wstąpić do [kogoś]"),
    @("verzweifeln an D", "zwątpić w [coś]"),
    @("vorbeireden an A", "mówić o czymś innym niż rozmówca, nie rozumieć się"),
    @("Anteil nehmen an D", "wziąć w czymś udział"),
    @("aus-rutschen auf D", "poślizgnąć się na [czymś]"),
    @("gut auskommen mit D", "wytyrzymywać z [kimś]"),
    @("jdm Bescheid geben über A", "poinformować kogoś/dać komuś znać o [czymś]"),
    @("sich bescheiden mit D", "zadowolić się [czymś]"),
    @("jdn überreden zu D", "namówić/przekonać kogoś do [czegoś]"),
    @("stolpern über A", "potknąć się o [coś]"),
    @("jdn unterstützen bei D", "pomagać komuś/wspierać kogoś w [czymś]"),
    @("sich verschanzen hinter D", "ukrywać się za [czymś]"),
    @("begeistert sein von D", "być zachwyconym [czymś]"),
)

$startRow = 216
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value() = $data[$i][0]
    $ws.Cells.Item($row, 3).Value() = $data[$i][1]
}

# Restore the selection state observed in the target workbook.
$ws.Range("A214").Select()
